$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        A = "Record"
        B = "Balanço Geral"
        C = "Saúde"
        D = "2025-04-10T12:34"
        E = "Neutro"
        F = "Alerta de chuva: fim de semana de tempo instável no Norte e Noroeste. Repórter *ao vivo*.Chuva pode chegar a qualquer momento. INMET emitiu alerta amarelo para cidades do Norte e Noroeste Fluminense. Ventos podem chegar a 60 km/h."
    },
    @{
        A = "Record"
        B = "Balanço Geral"
        C = "Saúde"
        D = "2025-04-10T12:37"
        E = "Negativo"
        F = "Hospital Plantadores de Cana paralisa alguns atendimentos a crianças. Repórter *ao vivo*. Funcionários há 3 meses sem receber. Reflexo no atendimento pediátrico de baixa complexidade, que acabou sendo suspenso na unidade. Na semana passada, os que atuam no Centro de Referência realizaram protesto em frente da unidade, que funciona anexo ao HPC. Em nota, Secretaria de Saúde mantém atendimento de emergência vermelha. Atendimento pediátrico também é feito na clínica da Criança. Repórter ressaltou que não esclarece sobre a falta de pagamento. "
    }
)

$startRow = 105
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("C$r").Value = $data.C
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
}
